$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the existing English-translation column (C4:C9) over to column D,
# making room for the new "Rewrite" notes in column C.
$ws.Range("C4:C9").Cut($ws.Range("D4:D9"))

# New rewrite note next to Lily's name (row 2).
$ws.Range("C2").Value = "Rewrite     -   Lily"

# New rewrite note next to Sina's name (row 30).
$ws.Range("C30").Value = "Sina"
